# Atualizando o arquivo XLSX
# Applies updated odds values to Sheet1, matching the upstream FlashScore refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "K2"  = 1.91
    "K4"  = 1.8
    "S4"  = 1.73
    "T4"  = 2.08
    "K5"  = 1.83
    "O5"  = 1.62
    "P5"  = 2.2
    "K6"  = 1.91
    "O6"  = 1.53
    "P6"  = 2.38
    "Q6"  = 2.7
    "R6"  = 1.44
    "M7"  = 1.14
    "N7"  = 5.5
    "O7"  = 1.62
    "M8"  = 1.04
    "O8"  = 1.25
    "G11" = 2.05
    "I11" = 4
    "K11" = 2.05
    "L11" = 4.33
    "U11" = 1.91
    "V11" = 1.8
    "W11" = 6.5
    "X11" = 9
    "AC11" = 8
    "AE11" = 15
    "AK11" = 34
    "AM11" = 351
    "AO11" = 12
    "AP11" = 23
    "AU11" = 8.5
    "AX11" = 21
    "BA11" = 101
    "BB11" = 251
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
